{"js": "// Apply targeted text replacements (date line + each division-problem cell)\n// as described by the diff. Every old string is unique within the\n// document, so a direct search-and-replace for each pair is unambiguous.\nconst replacements = [\n  [\"2024-07-27 Saturday\", \"2024-07-28 Sunday\"],\n  [\"289\u00f78=36, 1\", \"871\u00f76=145, 1\"],\n  [\"146\u00f77=20, 6\", \"707\u00f76=117, 5\"],\n  [\"992\u00f79=110, 2\", \"874\u00f72=437, 0\"],\n  [\"788\u00f77=112, 4\", \"450\u00f77=64, 2\"],\n  [\"595\u00f79=66, 1\", \"343\u00f75=68, 3\"],\n  [\"220\u00f73=73, 1\", \"401\u00f77=57, 2\"],\n  [\"943\u00f79=104, 7\", \"568\u00f72=284, 0\"],\n  [\"572\u00f73=190, 2\", \"548\u00f76=91, 2\"],\n  [\"137\u00f73=45, 2\", \"145\u00f72=72, 1\"],\n  [\"335\u00f75=67, 0\", \"198\u00f79=22, 0\"],\n  [\"666\u00f74=166, 2\", \"175\u00f75=35, 0\"],\n  [\"724\u00f74=181, 0\", \"317\u00f76=52, 5\"],\n  [\"153\u00f77=21, 6\", \"313\u00f77=44, 5\"],\n  [\"304\u00f75=60, 4\", \"890\u00f73=296, 2\"],\n  [\"326\u00f72=163, 0\", \"900\u00f73=300, 0\"],\n  [\"691\u00f74=172, 3\", \"894\u00f78=111, 6\"],\n  [\"364\u00f76=60, 4\", \"209\u00f75=41, 4\"],\n  [\"993\u00f77=141, 6\", \"518\u00f79=57, 5\"],\n  [\"537\u00f74=134, 1\", \"899\u00f75=179, 4\"],\n  [\"915\u00f76=152, 3\", \"194\u00f77=27, 5\"],\n  [\"114\u00f78=14, 2\", \"429\u00f78=53, 5\"],\n  [\"910\u00f78=113, 6\", \"610\u00f75=122, 0\"],\n  [\"639\u00f78=79, 7\", \"721\u00f79=80, 1\"],\n  [\"867\u00f79=96, 3\", \"577\u00f72=288, 1\"],\n  [\"211\u00f74=52, 3\", \"966\u00f75=193, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Apply targeted text replacements (date line + each division-problem\n# cell) as described by the diff. Every \"old\" string occurs exactly once\n# in the document, so Find/Replace (wdReplaceAll) for each pair is\n# unambiguous and leaves all other content untouched.\n$replacements = @(\n  @(\"2024-07-27 Saturday\", \"2024-07-28 Sunday\"),\n  @(\"289\u00f78=36, 1\", \"871\u00f76=145, 1\"),\n  @(\"146\u00f77=20, 6\", \"707\u00f76=117, 5\"),\n  @(\"992\u00f79=110, 2\", \"874\u00f72=437, 0\"),\n  @(\"788\u00f77=112, 4\", \"450\u00f77=64, 2\"),\n  @(\"595\u00f79=66, 1\", \"343\u00f75=68, 3\"),\n  @(\"220\u00f73=73, 1\", \"401\u00f77=57, 2\"),\n  @(\"943\u00f79=104, 7\", \"568\u00f72=284, 0\"),\n  @(\"572\u00f73=190, 2\", \"548\u00f76=91, 2\"),\n  @(\"137\u00f73=45, 2\", \"145\u00f72=72, 1\"),\n  @(\"335\u00f75=67, 0\", \"198\u00f79=22, 0\"),\n  @(\"666\u00f74=166, 2\", \"175\u00f75=35, 0\"),\n  @(\"724\u00f74=181, 0\", \"317\u00f76=52, 5\"),\n  @(\"153\u00f77=21, 6\", \"313\u00f77=44, 5\"),\n  @(\"304\u00f75=60, 4\", \"890\u00f73=296, 2\"),\n  @(\"326\u00f72=163, 0\", \"900\u00f73=300, 0\"),\n  @(\"691\u00f74=172, 3\", \"894\u00f78=111, 6\"),\n  @(\"364\u00f76=60, 4\", \"209\u00f75=41, 4\"),\n  @(\"993\u00f77=141, 6\", \"518\u00f79=57, 5\"),\n  @(\"537\u00f74=134, 1\", \"899\u00f75=179, 4\"),\n  @(\"915\u00f76=152, 3\", \"194\u00f77=27, 5\"),\n  @(\"114\u00f78=14, 2\", \"429\u00f78=53, 5\"),\n  @(\"910\u00f78=113, 6\", \"610\u00f75=122, 0\"),\n  @(\"639\u00f78=79, 7\", \"721\u00f79=80, 1\"),\n  @(\"867\u00f79=96, 3\", \"577\u00f72=288, 1\"),\n  @(\"211\u00f74=52, 3\", \"966\u00f75=193, 1\"),\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $range = $d.Content\n  $find = $range.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Forward = $true\n  $find.Wrap = 1            # wdFindContinue\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n\n  $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n  if (-not $found) {\n    throw \"Text not found: $oldText\"\n  }\n}\n"}
